$d = $word.ActiveDocument

# The five lines that already exist as the last five paragraphs of the
# document need to be duplicated and the duplicate inserted right after
# the leading empty paragraph (i.e. before the existing copy).
$texts = @(
    "D_Data_TCNo_72_73 _Approval :SUCESS: Approval Tab present.",
    "D_Data_TCNo_72_73 _Approval :SUCESS: Approver 1 name matched with expected.",
    "D_Data_TCNo_72_73 _Approval :SUCESS: Approver 2 drop down present.",
    "D_Data_TCNo_72_73 _Approval :SUCESS: Approver 3 drop down present.",
    "D_Data_TCNo_71 _SUCESS: Approval Send successfully."
)

# Anchor on what is currently the second paragraph (the first real text
# paragraph, right after the leading blank one) and insert a new blank
# paragraph before it for every line we need to add.
$anchor = $d.Paragraphs(2)
foreach ($t in $texts) {
    $anchor.Range.InsertParagraphBefore()
}

# Now fill in the freshly inserted (still empty) paragraphs, in order,
# with the corresponding text.
for ($i = 0; $i -lt $texts.Length; $i++) {
    $p = $d.Paragraphs(2 + $i)
    $p.Range.Text = $texts[$i]
}
